$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of an existing header cell (bold, centered, bordered)
# onto the three new header cells so they reuse the same cell style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Header labels for the new columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (same W/L/T for every player row, rows 2-45).
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 92
    $ws.Cells.Item($r, 31).Value = 70
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Output "done"
